$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 update: Inventory Quantity changes from 100 to 88 ---
$ws.Range("D2").Value = 88

# --- Row 3: new item "Harry Potter" ---
$ws.Range("B3").Value = "WL0047968"
$ws.Range("C3").Value = "Harry Potter"
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 5
$ws.Range("G3").Value = 10000
$ws.Range("H3").Value = 20000
$ws.Range("I3").Formula = "=H3*E3"
$ws.Range("J3").Formula = "=D3*G3"
$ws.Range("K3").Formula = "=I3-J3"

# --- Rows 4:25 get the same Total/Invested/Clean-profit formula pattern ---
for ($r = 4; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Formula  = "=H$r*E$r"
    $ws.Cells.Item($r, 10).Formula = "=D$r*G$r"
    $ws.Cells.Item($r, 11).Formula = "=I$r-J$r"
}

# --- Number formats (re-use the workbook's existing currency format, numFmtId 164) ---
$fmt164 = '_-[$₩-412]* #,##0.00_-;\-[$₩-412]* #,##0.00_-;_-[$₩-412]* "-"??_-;_-@_-'

# Header cells G1:K1 take on the bordered/filled header look + the currency format
$ws.Range("G1:K1").NumberFormat = $fmt164

# G3:H11 and I3:I25 take the centered currency style (matches existing G2:I2 style)
$r1 = $ws.Range("G3:H11")
$r1.NumberFormat = $fmt164
$r1.HorizontalAlignment = -4108

$r1b = $ws.Range("I3:I25")
$r1b.NumberFormat = $fmt164
$r1b.HorizontalAlignment = -4108

# J3:K25 take the plain currency style (matches existing J2:K2 style)
$r2 = $ws.Range("J3:K25")
$r2.NumberFormat = $fmt164

# --- Selection moves to H7, matching the saved workbook state ---
$ws.Range("H7").Select()

Write-Output "done"
